$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34; this shifts existing rows 34-154 down to 35-155.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record.
$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value = 44487
$ws.Range("D34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100102
$ws.Range("H34").Value = "Cítricos"
$ws.Range("I34").Value = 100102006
$ws.Range("J34").Value = "Pomelo"
$ws.Range("K34").Value = "Start Ruby"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 60
$ws.Range("N34").Value = 11000
$ws.Range("O34").Value = 12000
$ws.Range("P34").Value = 11500
$ws.Range("Q34").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R34").Value = "Región de O'Higgins"
$ws.Range("S34").Value = 821
$ws.Range("T34").Value = 14
